$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = "2022-12-08 17:37:14.919646"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "369"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = 6.99

# Row 12
$ws.Range("A12").Value = "2024-01-18 12:36:01.050333"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "789"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = 11.77
